# Insert a new "LP solver (linprog or gurobi)" = "gurobi" row into the
# "general" sheet, right before the existing "Number of exp. conditions"
# row (currently row 5), shifting everything below it down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# Insert a new row at position 5; existing rows 5.. shift down to 6..
$ws.Rows.Item(5).Insert()

$ws.Cells.Item(5, 1).Value = "LP solver (linprog or gurobi)"
$ws.Cells.Item(5, 2).Value = "gurobi"

# Make sure this sheet/cell is the active selection, matching the source
# workbook which now opens on the "general" tab.
$ws.Activate()
$ws.Range("A5").Select()
